$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers in E1 and F1, matching the formatting of the existing header row
$ws.Range("E1").Value = "Execution Time (ms)"
$ws.Range("F1").Value = "Memory Usage (B)"

# Copy the header style (bold font, border, centered/top alignment) from D1 to E1:F1
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new data columns
$ws.Range("E2").Value = 9.65699998778291
$ws.Range("F2").Value = 4096

$ws.Range("E3").Value = 9.225599991623312
$ws.Range("F3").Value = 4096

$ws.Range("E4").Value = 26.89470001496375
$ws.Range("F4").Value = 8192

$ws.Range("E5").Value = 6.516599998576567
$ws.Range("F5").Value = 0

$ws.Range("E6").Value = 1.92470001638867
$ws.Range("F6").Value = 0
